# A new record row is inserted into the dataset at row 131 (weekly price
# update for Berenjena / Vega Monumental Concepción). This pushes every
# existing row from 131..191 down by one (132..192) and grows the sheet's
# used range from A1:R191 to A1:R192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 131; Excel shifts rows 131-191 down
# to 132-192 and extends the sheet dimension automatically.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new record's data.
$ws.Range("A131").Value = 11
$ws.Range("B131").Value = "Vega Monumental Concepción"
$ws.Range("C131").Value = "Bíobío"
$ws.Range("D131").Value = Get-Date -Year 2023 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("E131").Value = 8
$ws.Range("F131").Value = 100112001
$ws.Range("G131").Value = "Berenjena"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 50
$ws.Range("K131").Value = 9000
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = 9000
$ws.Range("N131").Value = "`$/caja 50 unidades"
$ws.Range("O131").Value = "Región de Arica y Parinacota"
$ws.Range("P131").Value = 180
$ws.Range("Q131").Value = 50
$ws.Range("R131").Value = "Hortaliza"
